# This workbook tracks weekly "Zapallo italiano" price records.
# A new week of data (Huracán, Primera/Segunda) is inserted at the top of the
# data block (row 152 onward), pushing all the older weekly records down by
# two rows (the table keeps growing, oldest rows remain at the bottom, two
# brand-new rows appear for the new week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the first historical data row that
# needs to move down (row 152). Everything from 152 onward shifts to 154
# onward; the sheet grows from 270 to 272 used rows (dimension auto-updates).
$ws.Range("A152:A153").EntireRow.Insert()

# --- New row 152: Huracán / Primera, newest week ---
$ws.Cells.Item(152, 1).Value = 1
$ws.Cells.Item(152, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(152, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(152, 4).Value = 44603
$ws.Cells.Item(152, 5).Value = 15
$ws.Cells.Item(152, 6).Value = 100112032
$ws.Cells.Item(152, 7).Value = "Zapallo italiano"
$ws.Cells.Item(152, 8).Value = "Huracán"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 160
$ws.Cells.Item(152, 11).Value = 4500
$ws.Cells.Item(152, 12).Value = 5000
$ws.Cells.Item(152, 13).Value = 4750
$ws.Cells.Item(152, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(152, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(152, 16).Value = 68
$ws.Cells.Item(152, 17).Value = 70
$ws.Cells.Item(152, 18).Value = "Hortaliza"

# --- New row 153: Huracán / Segunda, newest week ---
$ws.Cells.Item(153, 1).Value = 1
$ws.Cells.Item(153, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(153, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(153, 4).Value = 44603
$ws.Cells.Item(153, 5).Value = 15
$ws.Cells.Item(153, 6).Value = 100112032
$ws.Cells.Item(153, 7).Value = "Zapallo italiano"
$ws.Cells.Item(153, 8).Value = "Huracán"
$ws.Cells.Item(153, 9).Value = "Segunda"
$ws.Cells.Item(153, 10).Value = 160
$ws.Cells.Item(153, 11).Value = 4000
$ws.Cells.Item(153, 12).Value = 4500
$ws.Cells.Item(153, 13).Value = 4250
$ws.Cells.Item(153, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(153, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(153, 16).Value = 42
$ws.Cells.Item(153, 17).Value = 100
$ws.Cells.Item(153, 18).Value = "Hortaliza"
